$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 94

$ws.Cells.Item($newRow, 1).Value = "19-12-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹13,418 per gram for 24 karat gold, ₹12,300 per gram for 22 karat gold and ₹10,064 per gram for 18 karat gold (also called 999 gold)."
